$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the account-statement ("Estado de Cuenta") rows 16-19:
#  - "Periodo Mora" (col E) periods are now listed most-recent-first
#  - "Valor Mora" (col F) keeps following its own period (rows reordered)
#  - "Salario Basico" (col G) is unified to 1,000,000 for every worker/period
$ws.Range("E16").Value = "2206"
$ws.Range("F16").Value = 29333
$ws.Range("G16").Value = 1000000

$ws.Range("E17").Value = "2112"
$ws.Range("F17").Value = 36341
$ws.Range("G17").Value = 1000000

$ws.Range("E18").Value = "2009"
$ws.Range("F18").Value = 35112
$ws.Range("G18").Value = 1000000

$ws.Range("E19").Value = "1908"
$ws.Range("F19").Value = 33125
$ws.Range("G19").Value = 1000000
